# Scheduled-runner refresh of market-board derived profit figures across
# the Sheets workbook. Only numeric "snapshot" cells (currentAveragePrice*,
# LevePrice*, LeveProfit*) are refreshed row-by-row; no formulas exist in
# this workbook, so each cell is written as a literal value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 13228.429
$ws.Range("I9").Value = 30356.334
$ws.Range("K9").Value = 30356.334
$ws.Range("M9").Value = -30187.334

$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 300
$ws.Range("K18").Value = 300
$ws.Range("M18").Value = -16

$ws.Range("H40").Value = 1218.0605
$ws.Range("I40").Value = 1194.4445
$ws.Range("K40").Value = 1194.4445
$ws.Range("M40").Value = -1019.4445

$ws.Range("H70").Value = 2298.6
$ws.Range("J70").Value = 2331.8333
$ws.Range("L70").Value = 6995.499899999999
$ws.Range("N70").Value = -7535.499899999999

$ws.Range("H73").Value = 2298.6
$ws.Range("J73").Value = 2331.8333
$ws.Range("L73").Value = 6995.499899999999
$ws.Range("N73").Value = -8867.499899999999

$ws.Range("H80").Value = 1648.3334
$ws.Range("J80").Value = 1599.1
$ws.Range("L80").Value = 4797.299999999999
$ws.Range("N80").Value = -6793.299999999999

$ws.Range("H83").Value = 1648.3334
$ws.Range("J83").Value = 1599.1
$ws.Range("L83").Value = 14391.9
$ws.Range("N83").Value = -24375.9

$ws.Range("H106").Value = 1572.4
$ws.Range("I106").Value = 990.5
$ws.Range("K106").Value = 990.5
$ws.Range("M106").Value = -359.5

$ws.Range("H116").Value = 10530.917
$ws.Range("I116").Value = 10636.5
$ws.Range("K116").Value = 10636.5
$ws.Range("M116").Value = -7194.5

$ws.Range("H121").Value = 1134.6
$ws.Range("J121").Value = 1249.5
$ws.Range("L121").Value = 3748.5
$ws.Range("N121").Value = -7242.5

$ws.Range("H132").Value = 9863745
$ws.Range("I132").Value = 10132343
$ws.Range("K132").Value = 30397029
$ws.Range("M132").Value = -30394499

$ws.Range("H137").Value = 1625.2354
$ws.Range("I137").Value = 1287
$ws.Range("J137").Value = 2724.5
$ws.Range("K137").Value = 3861
$ws.Range("L137").Value = 8173.5
$ws.Range("M137").Value = -1311
$ws.Range("N137").Value = -13273.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1336.52
$ws.Range("I97").Value = 1372.8182
$ws.Range("J97").Value = 1070.3334
$ws.Range("K97").Value = 1372.8182
$ws.Range("L97").Value = 1070.3334
$ws.Range("M97").Value = -876.8181999999999
$ws.Range("N97").Value = -2062.3334

$ws.Range("H110").Value = 9251.322
$ws.Range("I110").Value = 13502.066
$ws.Range("K110").Value = 13502.066
$ws.Range("M110").Value = -11457.066

$ws.Range("H132").Value = 4787.0713
$ws.Range("I132").Value = 5556.3335
$ws.Range("J132").Value = 3402.4
$ws.Range("K132").Value = 16669.0005
$ws.Range("L132").Value = 10207.2
$ws.Range("M132").Value = -14139.0005
$ws.Range("N132").Value = -15267.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1702.1351
$ws.Range("I134").Value = 1289.1936
$ws.Range("J134").Value = 3835.6667
$ws.Range("K134").Value = 3867.5808
$ws.Range("L134").Value = 11507.0001
$ws.Range("M134").Value = -1332.5808
$ws.Range("N134").Value = -16577.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 173.76471
$ws.Range("J7").Value = 370.83334
$ws.Range("L7").Value = 370.83334
$ws.Range("N7").Value = -596.83334

$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1574

$ws.Range("H22").Value = 292.625
$ws.Range("J22").Value = 268.2
$ws.Range("L22").Value = 268.2
$ws.Range("N22").Value = -968.2

$ws.Range("H32").Value = 3050.8572
$ws.Range("I32").Value = 3050.8572
$ws.Range("K32").Value = 3050.8572
$ws.Range("M32").Value = -2734.8572

$ws.Range("H99").Value = 3860
$ws.Range("I99").Value = 3794.1667
$ws.Range("J99").Value = 4156.25
$ws.Range("K99").Value = 3794.1667
$ws.Range("L99").Value = 4156.25
$ws.Range("M99").Value = -2296.1667
$ws.Range("N99").Value = -7152.25

$ws.Range("H105").Value = 2830
$ws.Range("I105").Value = 3202.5
$ws.Range("K105").Value = 3202.5
$ws.Range("M105").Value = -1455.5

$ws.Range("H107").Value = 986.6923
$ws.Range("I107").Value = 938.8182
$ws.Range("K107").Value = 938.8182
$ws.Range("M107").Value = 981.1818

$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("N113").Value = -5340

$ws.Range("H126").Value = 3860
$ws.Range("I126").Value = 3794.1667
$ws.Range("J126").Value = 4156.25
$ws.Range("K126").Value = 11382.5001
$ws.Range("L126").Value = 12468.75
$ws.Range("M126").Value = -8912.500100000001
$ws.Range("N126").Value = -17408.75

$ws.Range("H134").Value = 15768.72
$ws.Range("I134").Value = 7421.8945
$ws.Range("J134").Value = 42200.332
$ws.Range("K134").Value = 22265.6835
$ws.Range("L134").Value = 126600.996
$ws.Range("M134").Value = -19730.6835
$ws.Range("N134").Value = -131670.996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 50000000
$ws.Range("J68").Value = 50000000
$ws.Range("L68").Value = 150000000
$ws.Range("N68").Value = -150001622

$ws.Range("H71").Value = 50000000
$ws.Range("J71").Value = 50000000
$ws.Range("L71").Value = 450000000
$ws.Range("N71").Value = -450008112

$ws.Range("H99").Value = 6739

$ws.Range("H132").Value = 1582.3529
$ws.Range("J132").Value = 2283.3333
$ws.Range("L132").Value = 20549.9997
$ws.Range("N132").Value = -25609.9997

$ws.Range("H139").Value = 5873.375
$ws.Range("I139").Value = 6141
$ws.Range("K139").Value = 18423
$ws.Range("M139").Value = -13283

$ws.Range("H140").Value = 21155.334
$ws.Range("I140").Value = 21155.334
$ws.Range("K140").Value = 63466.00199999999
$ws.Range("M140").Value = -58286.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 70000
$ws.Range("J39").Value = 70000
$ws.Range("L39").Value = 70000
$ws.Range("N39").Value = -71064

$ws.Range("H57").Value = 6183

$ws.Range("H122").Value = 1318.1111
$ws.Range("I122").Value = 1181.8572
$ws.Range("K122").Value = 3545.5716
$ws.Range("M122").Value = -1095.5716

$ws.Range("H132").Value = 230025.8
$ws.Range("I132").Value = 259035.97
$ws.Range("J132").Value = 3746.4
$ws.Range("K132").Value = 777107.91
$ws.Range("L132").Value = 11239.2
$ws.Range("M132").Value = -774577.91
$ws.Range("N132").Value = -16299.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H40").Value = 4584.25
$ws.Range("I40").Value = 3944.8333
$ws.Range("K40").Value = 3944.8333
$ws.Range("M40").Value = -3808.8333

$ws.Range("H93").Value = 588042.3
$ws.Range("I93").Value = 795671.8
$ws.Range("K93").Value = 795671.8
$ws.Range("M93").Value = -794423.8

$ws.Range("H136").Value = 4070.75
$ws.Range("I136").Value = 3781.862
$ws.Range("K136").Value = 11345.586
$ws.Range("M136").Value = -8795.585999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13500.2
$ws.Range("J45").Value = 14126
$ws.Range("L45").Value = 14126
$ws.Range("N45").Value = -15108

$ws.Range("H96").Value = 2242
$ws.Range("I96").Value = 1875
$ws.Range("J96").Value = 2364.3333
$ws.Range("K96").Value = 1875
$ws.Range("L96").Value = 2364.3333
$ws.Range("M96").Value = -502
$ws.Range("N96").Value = -5110.3333

$ws.Range("H122").Value = 2423.3044
$ws.Range("I122").Value = 2372.5789
$ws.Range("J122").Value = 2664.25
$ws.Range("K122").Value = 7117.736699999999
$ws.Range("L122").Value = 7992.75
$ws.Range("M122").Value = -4667.736699999999
$ws.Range("N122").Value = -12892.75

$ws.Range("H126").Value = 7377.2
$ws.Range("I126").Value = 8226.25
$ws.Range("K126").Value = 24678.75
$ws.Range("M126").Value = -22208.75

$ws.Range("H132").Value = 2717.6843
$ws.Range("I132").Value = 2646.7222
$ws.Range("J132").Value = 3995
$ws.Range("K132").Value = 7940.1666
$ws.Range("L132").Value = 11985
$ws.Range("M132").Value = -5410.1666
$ws.Range("N132").Value = -17045
